$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - update related_works (column AF)
$ws.Range("AF6").Value = 'c("https://openalex.org/W4238752995", "https://openalex.org/W2393977738", "https://openalex.org/W2088706544", "https://openalex.org/W4207046597", "https://openalex.org/W4285669903", "https://openalex.org/W2051801047", "https://openalex.org/W1994583257", "https://openalex.org/W1867320779", "https://openalex.org/W2142353978", "https://openalex.org/W1973582222")'

# Row 10 - update related_works (column AF)
$ws.Range("AF10").Value = 'c("https://openalex.org/W2938290407", "https://openalex.org/W4382753601", "https://openalex.org/W2588901189", "https://openalex.org/W2767430559", "https://openalex.org/W2751366772", "https://openalex.org/W4281722455", "https://openalex.org/W3007362983", "https://openalex.org/W4295248011", "https://openalex.org/W2772482411", "https://openalex.org/W2468080474")'

# Row 13 - fill in source/journal info (columns F, G, H, I)
$ws.Range("F13").Value = "Chest"
$ws.Range("G13").Value = "https://openalex.org/S76900504"
$ws.Range("H13").Value = "Elsevier BV"
$ws.Range("I13").Value = "0012-3692"

# Row 13 - fix any_repository_has_fulltext (column V)
# Leading apostrophe forces Excel to store this as text "FALSE" rather than a Boolean,
# matching the original inlineStr text representation used in this sheet.
$ws.Range("V13").Value = "'FALSE"
